$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert new row 2 with "(percent)" labels in three languages ---
$ws.Rows("2:2").Insert()
$ws.Rows("2:2").Clear()

$ws.Range("A2").Value = "(пайыз менен)"
$ws.Range("B2").Value = "(в процентах)"
$ws.Range("C2").Value = "(in percent)"

$hdrRng = $ws.Range("A2:C2")
$hdrRng.Font.Name = "Times New Roman"
$hdrRng.Font.Italic = $true
$hdrRng.Font.Size = 8
$hdrRng.HorizontalAlignment = -4108
$hdrRng.VerticalAlignment = -4108
$hdrRng.WrapText = $true

$ws.Rows(2).RowHeight = 14.25
$ws.Rows(3).RowHeight = 14.25

# --- 2. Add new column O (year 2023) mirroring column N's formatting ---
$ws.Range("N4").Copy()
$ws.Range("O4").PasteSpecial(-4122)
$ws.Range("O4").Value = 2023

$ws.Range("N5").Copy()
$ws.Range("O5").PasteSpecial(-4122)
$ws.Range("O5").Value = 7.9591147916539313

$ws.Range("N6").Copy()
$ws.Range("O6").PasteSpecial(-4122)
$ws.Range("O6").Value = 4.1262815690193904

$ws.Range("N7").Copy()
$ws.Range("O7").PasteSpecial(-4122)
$ws.Range("O7").Value = 11.553674062171684

$ws.Range("N8").Copy()
$ws.Range("O8").PasteSpecial(-4122)
$ws.Range("O8").Value = 26.7840134279745

$ws.Range("N9").Copy()
$ws.Range("O9").PasteSpecial(-4122)
$ws.Range("O9").Value = 31.703252552185106

$ws.Range("N10").Copy()
$ws.Range("O10").PasteSpecial(-4122)
$ws.Range("O10").Value = 22.127282549972989

$ws.Range("N11").Copy()
$ws.Range("O11").PasteSpecial(-4122)
$ws.Range("O11").Value = 25.785751793343863

$ws.Range("N12").Copy()
$ws.Range("O12").PasteSpecial(-4122)
$ws.Range("O12").Value = 27.265979822798002

$ws.Range("N13").Copy()
$ws.Range("O13").PasteSpecial(-4122)
$ws.Range("O13").Value = 24.322552749117975

$excel.CutCopyMode = $false

# --- 3. Clear the stale selection so the saved sheetView has none ---
$ws.Range("A1").Select()

Write-Output "edit complete"
